$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Shift the "Score" row (was row 9) up to row 8, freeing row 9 for a new
#     "Probability" row, and row 10 becomes a red "final answer" cell that
#     either shows the probability or a "Can not be determined" warning. ---

# 1) Row 8: add "Score" label in C8 (copy format from the old C9 "Score" label)
#    and move the scoring formula from D9 into D8.
$ws.Range("C9").Copy($ws.Range("C8"))
$ws.Range("C8").Value = "Score"

$ws.Range("D9").Copy($ws.Range("D8"))
$ws.Range("D8").Formula = '=$B$2*D2+$B$3*D3+$B$4*D4+$B$5*D5+(D6-D7)*$B$8+$B$9'

# 2) Row 9: relabel C9 to "Probability" (matching the style used by the other
#    hidden/white helper cells such as A10/B10/B11) and put the probability
#    formula (now referencing D8 instead of D9) into D9.
$ws.Range("B11").Copy($ws.Range("C9"))
$ws.Range("C9").Value = "Probability"

$ws.Range("B11").Copy($ws.Range("D9"))
$ws.Range("D9").Formula = '=$B$11*EXP(D8)/(1+$B$11*EXP(D8))'

# 3) Row 10: D10 becomes the visible result - a warning text if the
#    probability is too low to be meaningful, otherwise the probability
#    itself. Keep its red, boxed styling but resize to match the rest.
$ws.Range("D10").Formula = '=IF(D9<0.5, "Can not be determined", D9)'
$ws.Range("D10").Font.Size = 11

# --- Cosmetics that accompanied the edit ---

# New column for the (now wider) "Can not be determined" text.
$ws.Columns(4).AutoFit()

# Cursor / selection left on C11 when the file was saved.
$ws.Range("C11").Select()
